# Update metadata in the sites
# Adds "instrument"/"canopy" columns and per-station "Metadata" link columns
# to the site list table (Sheet1), matching the upstream commit that added
# ICOS metadata / instrument-canopy-height info for the stations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters here: it drives the order new entries are appended to the
# shared-strings table, so we replay them in the same sequence the author
# apparently entered them in (IT-Tor metadata link first, then the new
# header cells, then the BE-Lon / Lonzee links, then DE-HoH, then the new
# "canopy"/"instrument" headers).

# Row 5 (Torgnon / IT-Tor): Metadata link
$ws.Range("J5").Value = "https://meta.icos-cp.eu/resources/stations/ES_IT-Tor"

# New header cells on row 1
$ws.Range("J1").Value = "Metadata"

# Row 3 (Lonzee / BE-Lon): two metadata-ish links
$ws.Range("J3").Value = "https://www.icos-belgium.be/ESLonzee.php"
$ws.Range("K3").Value = "https://meta.icos-cp.eu/resources/stations/ES_BE-Lon"

# Row 4 (Hohes Holz / DE-HoH): Metadata link
$ws.Range("J4").Value = "https://meta.icos-cp.eu/resources/stations/ES_DE-HoH"

# Remaining new headers on row 1
$ws.Range("I1").Value = "canopy"
$ws.Range("H1").Value = "instrument"

# Row 4 instrument/canopy height values
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = 33

# Matches the workbook's saved selection (cell below the new columns)
$ws.Range("I7").Select()
